# Handback status report regeneration: bump the "generate date" /
# handoff- and handback-datetime timestamps recorded for this run.
#
#   Overview!G2        (Latest HO Xliff Generate Date)        01:01:20 -> 01:02:31
#   zh-cn!H2            (Correspond Handoff Datetime)          01:01:14 -> 01:02:25
#   zh-cn!K2            (Correspond Handback Datetime)         01:01:50 -> 01:02:42
#   de-de!H2            (Correspond Handoff Datetime)          01:01:20 -> 01:02:31
#   de-de!K2            (Correspond Handback Datetime)         01:01:57 -> 01:02:49

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-25 01:02:31"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-25 01:02:25"
$zhcn.Range("K2").Value = "2016-08-25 01:02:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-25 01:02:31"
$dede.Range("K2").Value = "2016-08-25 01:02:49"
